$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lương")

$ws.Range("B2").Value = 8
$ws.Range("B3").Value = 571428.5714285715
$ws.Range("B13").Value = 857142.8571428573
$ws.Range("B23").Value = 857142.8571428573
$ws.Range("B31").Value = 571428.5714285715
$ws.Range("B32").Value = 857142.8571428573
$ws.Range("B33").Value = 857142.8571428573
$ws.Range("B34").Value = 2285714.285714286
